# Update "想去人数" (column F) counts that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 41
$ws1.Range("F3").Value = 110
$ws1.Range("F4").Value = 1543
$ws1.Range("F5").Value = 235
$ws1.Range("F7").Value = 751
$ws1.Range("F8").Value = 10052
$ws1.Range("F9").Value = 172
$ws1.Range("F10").Value = 128
$ws1.Range("F12").Value = 190
$ws1.Range("F13").Value = 381
$ws1.Range("F14").Value = 6965
$ws1.Range("F15").Value = 1091
$ws1.Range("F16").Value = 651
$ws1.Range("F17").Value = 55

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 7

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 41
$ws4.Range("F3").Value = 110
$ws4.Range("F4").Value = 1543
$ws4.Range("F5").Value = 235
$ws4.Range("F6").Value = 7
$ws4.Range("F8").Value = 751
$ws4.Range("F11").Value = 10052
$ws4.Range("F12").Value = 172
$ws4.Range("F13").Value = 128
$ws4.Range("F15").Value = 190
$ws4.Range("F16").Value = 381
$ws4.Range("F17").Value = 6965
$ws4.Range("F18").Value = 1091
$ws4.Range("F19").Value = 651
$ws4.Range("F20").Value = 55

$wb.Save()
